$wb = $excel.ActiveWorkbook

# Update the "想去人数" (number of people interested) counts for two events.
# These values live both on the "展览" sheet and the aggregated "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 102
    $ws.Range("F10").Value = 404
}
